$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $escaped = $val -replace '"', '""'
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy() | Out-Null
    $rng.PasteSpecial(-4163) | Out-Null
}

Set-TextValue $ws.Range("D2") '37.101.17'
Set-TextValue $ws.Range("E2") '  -0.16%  '
Set-TextValue $ws.Range("D3") '2.068.20'
Set-TextValue $ws.Range("E3") '  -0.84%  '
Set-TextValue $ws.Range("E4") '  +0.15%  '
Set-TextValue $ws.Range("D5") '252.95'
Set-TextValue $ws.Range("E5") '  +1.11%  '
Set-TextValue $ws.Range("D6") '0.673'
Set-TextValue $ws.Range("E6") '  +1.80%  '
Set-TextValue $ws.Range("D7") '59.48'
Set-TextValue $ws.Range("E7") '  +10.55%  '
Set-TextValue $ws.Range("E8") '  -0.06%  '
Set-TextValue $ws.Range("D9") '61.49'
Set-TextValue $ws.Range("E9") '  -0.31%  '
Set-TextValue $ws.Range("D10") '0.387'
Set-TextValue $ws.Range("E10") '  +3.75%  '
Set-TextValue $ws.Range("E12") '  +2.51%  '
Set-TextValue $ws.Range("D13") '16.21'
Set-TextValue $ws.Range("E13") '  +7.93%  '
Set-TextValue $ws.Range("D14") '2.373.93'
Set-TextValue $ws.Range("E14") '  -0.69%  '
Set-TextValue $ws.Range("D15") '0.813'
Set-TextValue $ws.Range("E15") '  -1.82%  '
Set-TextValue $ws.Range("D16") '5.62'
Set-TextValue $ws.Range("E16") '  +9.10%  '
Set-TextValue $ws.Range("D17") '2.063.94'
Set-TextValue $ws.Range("E17") '  -1.12%  '
Set-TextValue $ws.Range("D18") '37.066.58'
Set-TextValue $ws.Range("E18") '  -0.24%  '
Set-TextValue $ws.Range("D19") '16.43'
Set-TextValue $ws.Range("E19") '  +12.79%  '
Set-TextValue $ws.Range("D20") '75.10'
Set-TextValue $ws.Range("E20") '  +3.45%  '
Set-TextValue $ws.Range("D21") '0.0₃0927'
Set-TextValue $ws.Range("E21") '  +9.84%  '
Set-TextValue $ws.Range("D22") '5.48'
Set-TextValue $ws.Range("E22") '  +5.66%  '
Set-TextValue $ws.Range("D23") '239.13'
Set-TextValue $ws.Range("E23") '  -0.43%  '
Set-TextValue $ws.Range("E24") '  -0.08%  '
Set-TextValue $ws.Range("E25") '  -2.08%  '
Set-TextValue $ws.Range("E26") '  +15.12%  '
Set-TextValue $ws.Range("D27") '169.22'
Set-TextValue $ws.Range("E27") '  -1.69%  '
Set-TextValue $ws.Range("D28") '9.33'
Set-TextValue $ws.Range("E28") '  +1.31%  '
Set-TextValue $ws.Range("D29") '20.35'
Set-TextValue $ws.Range("E29") '  -1.45%  '
Set-TextValue $ws.Range("E30") '  +3.05%  '
Set-TextValue $ws.Range("E31") '  +5.45%  '
Set-TextValue $ws.Range("D32") '4.79'
Set-TextValue $ws.Range("E32") '  +6.19%  '
Set-TextValue $ws.Range("D33") '0.0622'
Set-TextValue $ws.Range("E33") '  +1.21%  '
Set-TextValue $ws.Range("D34") '4.51'
Set-TextValue $ws.Range("E34") '  +8.87%  '
Set-TextValue $ws.Range("D35") '0.0910'
Set-TextValue $ws.Range("E35") '  +0.63%  '
Set-TextValue $ws.Range("E36") '  -0.04%  '
Set-TextValue $ws.Range("D37") '2.30'
Set-TextValue $ws.Range("E37") '  +3.59%  '
Set-TextValue $ws.Range("D38") '0.117'
Set-TextValue $ws.Range("E38") '  +24.47%  '
Set-TextValue $ws.Range("E39") '  -2.42%  '
Set-TextValue $ws.Range("D40") '1.37'
Set-TextValue $ws.Range("E40") '  +1.89%  '
Set-TextValue $ws.Range("D41") '17.95'
Set-TextValue $ws.Range("E41") '  -1.06%  '
Set-TextValue $ws.Range("D42") '0.0227'
Set-TextValue $ws.Range("E42") '  +0.52%  '
Set-TextValue $ws.Range("E43") '  +0.48%  '
Set-TextValue $ws.Range("D44") '98.48'
Set-TextValue $ws.Range("E44") '  +0.41%  '
Set-TextValue $ws.Range("E45") '  +2.17%  '
Set-TextValue $ws.Range("B46") 'THORChain'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D46") '4.63'
Set-TextValue $ws.Range("E46") '  +15.09%  '
Set-TextValue $ws.Range("B47") 'FTXToken'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range("D47") '4.09'
Set-TextValue $ws.Range("E47") '  -7.95%  '
Set-TextValue $ws.Range("D48") '2.49'
Set-TextValue $ws.Range("E48") '  +7.03%  '
Set-TextValue $ws.Range("D49") '1.302.19'
Set-TextValue $ws.Range("E49") '  -1.31%  '
Set-TextValue $ws.Range("D50") '2.92'
Set-TextValue $ws.Range("E50") '  -0.38%  '
Set-TextValue $ws.Range("D51") '6.94'
Set-TextValue $ws.Range("E51") '  -0.49%  '

$excel.CutCopyMode = 0
